# Re-run of a classifier benchmark: properties/timestamp updated and the
# metrics sheet regenerated with two changed rows plus six new rows.

$wb = $excel.ActiveWorkbook

# --- properties sheet: total execution time + timestamp ---
$wsProps = $wb.Worksheets.Item("properties")
$wsProps.Range("M2").Value = 39.59
$wsProps.Range("N2").Value = "Thu Jun 13 08:47:43 2024"

# --- metrics sheet ---
$wsMetrics = $wb.Worksheets.Item("metrics")

# Apply the header-style formatting (bold/border/centered, style index used by
# A2:A3) to the new index cells A4:A9 before filling in their values.
$wsMetrics.Range("A2:A3").Copy() | Out-Null
$wsMetrics.Range("A4:A9").PasteSpecial(-4122) | Out-Null

# Row 2: BernoulliNB(alpha=1) -> RandomForestClassifier()
$wsMetrics.Range("B2").Value = "RandomForestClassifier()"
$wsMetrics.Range("C2").Value = 0.7683253588516746
$wsMetrics.Range("D2").Value = 0.8329550586191924
$wsMetrics.Range("E2").Value = 0.02435033385987805
$wsMetrics.Range("F2").Value = 0.0146266557147158

# Row 3: mean_classifier -> KNeighborsClassifier(n_neighbors=20)
$wsMetrics.Range("B3").Value = "KNeighborsClassifier(n_neighbors=20)"
$wsMetrics.Range("C3").Value = 0.6843992557150451
$wsMetrics.Range("D3").Value = 0.8484459883244077
$wsMetrics.Range("E3").Value = 0.03452429237957232
$wsMetrics.Range("F3").Value = 0.006929826070852984

# Row 4: AdaBoostClassifier
$wsMetrics.Range("A4").Value = 2
$wsMetrics.Range("B4").Value = "AdaBoostClassifier(algorithm='SAMME', learning_rate=1)"
$wsMetrics.Range("C4").Value = 0.7622169059011165
$wsMetrics.Range("D4").Value = 0.8409233125874463
$wsMetrics.Range("E4").Value = 0.02999158647670358
$wsMetrics.Range("F4").Value = 0.01289061083875735
$wsMetrics.Range("G4").Value = "HRV_Heart_rate"
$wsMetrics.Range("H4").Value = "HRV_MeanNN"
$wsMetrics.Range("I4").Value = "HRV_SDNN"
$wsMetrics.Range("J4").Value = 700
$wsMetrics.Range("K4").Value = $true
$wsMetrics.Range("L4").Value = $false
$wsMetrics.Range("M4").Value = $false
$wsMetrics.Range("N4").Value = $false
$wsMetrics.Range("O4").Value = $false
$wsMetrics.Range("P4").Value = 60
$wsMetrics.Range("Q4").Value = "WESAD"
$wsMetrics.Range("R4").Value = $true
$wsMetrics.Range("S4").Value = $false
$wsMetrics.Range("T4").Value = $true

# Row 5: DecisionTreeClassifier
$wsMetrics.Range("A5").Value = 3
$wsMetrics.Range("B5").Value = "DecisionTreeClassifier(max_depth=2)"
$wsMetrics.Range("C5").Value = 0.6477777777777778
$wsMetrics.Range("D5").Value = 0.8119584599797366
$wsMetrics.Range("E5").Value = 0.03308113949591328
$wsMetrics.Range("F5").Value = 0.00583294529920401
$wsMetrics.Range("G5").Value = "HRV_Heart_rate"
$wsMetrics.Range("H5").Value = "HRV_MeanNN"
$wsMetrics.Range("I5").Value = "HRV_SDNN"
$wsMetrics.Range("J5").Value = 700
$wsMetrics.Range("K5").Value = $true
$wsMetrics.Range("L5").Value = $false
$wsMetrics.Range("M5").Value = $false
$wsMetrics.Range("N5").Value = $false
$wsMetrics.Range("O5").Value = $false
$wsMetrics.Range("P5").Value = 60
$wsMetrics.Range("Q5").Value = "WESAD"
$wsMetrics.Range("R5").Value = $true
$wsMetrics.Range("S5").Value = $false
$wsMetrics.Range("T5").Value = $true

# Row 6: SVC
$wsMetrics.Range("A6").Value = 4
$wsMetrics.Range("B6").Value = "SVC()"
$wsMetrics.Range("C6").Value = 0.7574401913875597
$wsMetrics.Range("D6").Value = 0.8492543783470835
$wsMetrics.Range("E6").Value = 0.03557855856474591
$wsMetrics.Range("F6").Value = 0.01858333540622863
$wsMetrics.Range("J6").Value = 700
$wsMetrics.Range("K6").Value = $true
$wsMetrics.Range("L6").Value = $false
$wsMetrics.Range("M6").Value = $false
$wsMetrics.Range("N6").Value = $false
$wsMetrics.Range("O6").Value = $false
$wsMetrics.Range("P6").Value = 60
$wsMetrics.Range("Q6").Value = "WESAD"
$wsMetrics.Range("R6").Value = $true
$wsMetrics.Range("S6").Value = $false
$wsMetrics.Range("T6").Value = $true

# Row 7: LinearDiscriminantAnalysis
$wsMetrics.Range("A7").Value = 5
$wsMetrics.Range("B7").Value = "LinearDiscriminantAnalysis()"
$wsMetrics.Range("C7").Value = 0.7442187163239795
$wsMetrics.Range("D7").Value = 0.8365865055241954
$wsMetrics.Range("E7").Value = 0.04157821370232159
$wsMetrics.Range("F7").Value = 0.02094162707532396
$wsMetrics.Range("G7").Value = "HRV_Heart_rate"
$wsMetrics.Range("H7").Value = "HRV_MeanNN"
$wsMetrics.Range("I7").Value = "HRV_SDNN"
$wsMetrics.Range("J7").Value = 700
$wsMetrics.Range("K7").Value = $true
$wsMetrics.Range("L7").Value = $false
$wsMetrics.Range("M7").Value = $false
$wsMetrics.Range("N7").Value = $false
$wsMetrics.Range("O7").Value = $false
$wsMetrics.Range("P7").Value = 60
$wsMetrics.Range("Q7").Value = "WESAD"
$wsMetrics.Range("R7").Value = $true
$wsMetrics.Range("S7").Value = $false
$wsMetrics.Range("T7").Value = $true

# Row 8: BernoulliNB(alpha=1)
$wsMetrics.Range("A8").Value = 6
$wsMetrics.Range("B8").Value = "BernoulliNB(alpha=1)"
$wsMetrics.Range("C8").Value = 0.7673821429084587
$wsMetrics.Range("D8").Value = 0.7543260963960052
$wsMetrics.Range("E8").Value = 0.01485837740238304
$wsMetrics.Range("F8").Value = 0.03754085420231147
$wsMetrics.Range("G8").Value = "HRV_Heart_rate"
$wsMetrics.Range("H8").Value = "HRV_MeanNN"
$wsMetrics.Range("I8").Value = "HRV_SDNN"
$wsMetrics.Range("J8").Value = 700
$wsMetrics.Range("K8").Value = $true
$wsMetrics.Range("L8").Value = $false
$wsMetrics.Range("M8").Value = $false
$wsMetrics.Range("N8").Value = $false
$wsMetrics.Range("O8").Value = $false
$wsMetrics.Range("P8").Value = 60
$wsMetrics.Range("Q8").Value = "WESAD"
$wsMetrics.Range("R8").Value = $true
$wsMetrics.Range("S8").Value = $false
$wsMetrics.Range("T8").Value = $true

# Row 9: mean_classifier
$wsMetrics.Range("A9").Value = 7
$wsMetrics.Range("B9").Value = "mean_classifier"
$wsMetrics.Range("C9").Value = 0.733108621266516
$wsMetrics.Range("D9").Value = 0.8249213999682953
$wsMetrics.Range("E9").Value = 0.03056607169735969
$wsMetrics.Range("F9").Value = 0.01676369351534203
$wsMetrics.Range("J9").Value = 700
$wsMetrics.Range("K9").Value = $true
$wsMetrics.Range("L9").Value = $false
$wsMetrics.Range("M9").Value = $false
$wsMetrics.Range("N9").Value = $false
$wsMetrics.Range("O9").Value = $false
$wsMetrics.Range("P9").Value = 60
$wsMetrics.Range("Q9").Value = "WESAD"
$wsMetrics.Range("R9").Value = $true
$wsMetrics.Range("S9").Value = $false
$wsMetrics.Range("T9").Value = $true
